# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the Belias Profits workbook (market price refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 138227.5
$ws.Range("J17").Value = 138227.5
$ws.Range("L17").Value = 414682.5
$ws.Range("N17").Value = -415018.5

$ws.Range("H82").Value = 1922.5
$ws.Range("I82").Value = 878.4
$ws.Range("J82").Value = 7143
$ws.Range("K82").Value = 2635.2
$ws.Range("L82").Value = 21429
$ws.Range("M82").Value = -2229.2
$ws.Range("N82").Value = -22241

$ws.Range("H85").Value = 1922.5
$ws.Range("I85").Value = 878.4
$ws.Range("J85").Value = 7143
$ws.Range("K85").Value = 2635.2
$ws.Range("L85").Value = 21429
$ws.Range("M85").Value = -1231.2
$ws.Range("N85").Value = -24237

$ws.Range("H113").Value = 4385.3145
$ws.Range("I113").Value = 4304.5835
$ws.Range("J113").Value = 4561.4546
$ws.Range("K113").Value = 4304.5835
$ws.Range("L113").Value = 4561.4546
$ws.Range("M113").Value = -1050.5835
$ws.Range("N113").Value = -11069.4546

$ws.Range("H116").Value = 3000
$ws.Range("J116").Value = 3000
$ws.Range("L116").Value = 3000
$ws.Range("N116").Value = -9884

$ws.Range("H132").Value = 3128.4707
$ws.Range("I132").Value = 1128.3043
$ws.Range("J132").Value = 21530
$ws.Range("K132").Value = 3384.9129
$ws.Range("L132").Value = 64590
$ws.Range("M132").Value = -854.9129000000003
$ws.Range("N132").Value = -69650

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 44289.145
$ws.Range("J9").Value = 32500
$ws.Range("L9").Value = 32500
$ws.Range("N9").Value = -32840

$ws.Range("H20").Value = 44289.145
$ws.Range("J20").Value = 32500
$ws.Range("L20").Value = 32500
$ws.Range("N20").Value = -33040

$ws.Range("H27").Value = 9166.666999999999
$ws.Range("J27").Value = 9166.666999999999
$ws.Range("L27").Value = 9166.666999999999
$ws.Range("N27").Value = -9534.666999999999

$ws.Range("H32").Value = 8347.15
$ws.Range("I32").Value = 5654.378
$ws.Range("J32").Value = 20614.223
$ws.Range("K32").Value = 5654.378
$ws.Range("L32").Value = 20614.223
$ws.Range("M32").Value = -5367.378
$ws.Range("N32").Value = -21188.223

$ws.Range("H61").Value = 1650.5
$ws.Range("I61").Value = 1650.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1650.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1438.5
$ws.Range("N61").ClearContents()

$ws.Range("H136").Value = 1650.5
$ws.Range("I136").Value = 1650.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4951.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2401.5
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 26008.666
$ws.Range("I32").Value = 8026
$ws.Range("J32").Value = 35000
$ws.Range("K32").Value = 8026
$ws.Range("L32").Value = 35000
$ws.Range("M32").Value = -7642
$ws.Range("N32").Value = -35768

$ws.Range("H33").Value = 27500

$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H46").Value = 20000
$ws.Range("J46").Value = 20000
$ws.Range("L46").Value = 20000
$ws.Range("N46").Value = -20596

$ws.Range("H86").Value = 2114.375
$ws.Range("I86").Value = 2056.625
$ws.Range("J86").Value = 2229.875
$ws.Range("K86").Value = 2056.625
$ws.Range("L86").Value = 2229.875
$ws.Range("M86").Value = -933.625
$ws.Range("N86").Value = -4475.875

$ws.Range("H89").Value = 2114.375
$ws.Range("I89").Value = 2056.625
$ws.Range("J89").Value = 2229.875
$ws.Range("K89").Value = 10283.125
$ws.Range("L89").Value = 11149.375
$ws.Range("M89").Value = -4667.125
$ws.Range("N89").Value = -22381.375

$ws.Range("H94").Value = 1239.7587
$ws.Range("I94").Value = 957.2083
$ws.Range("J94").Value = 2596
$ws.Range("K94").Value = 957.2083
$ws.Range("L94").Value = 2596
$ws.Range("M94").Value = -506.2083
$ws.Range("N94").Value = -3498

$ws.Range("H134").Value = 6633.7407
$ws.Range("I134").Value = 6559.0557
$ws.Range("K134").Value = 19677.1671
$ws.Range("M134").Value = -17142.1671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 9925
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 9925
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 9925
$ws.Range("N36").Value = -10701
$ws.Range("M36").ClearContents()

$ws.Range("H40").Value = 9925
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 9925
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 9925
$ws.Range("N40").Value = -10245
$ws.Range("M40").ClearContents()

$ws.Range("H132").Value = 1614.6038
$ws.Range("I132").Value = 1259.5897
$ws.Range("J132").Value = 2603.5715
$ws.Range("K132").Value = 3778.7691
$ws.Range("L132").Value = 7810.7145
$ws.Range("M132").Value = -1248.7691
$ws.Range("N132").Value = -12870.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 942.9167
$ws.Range("I4").Value = 41.25
$ws.Range("K4").Value = 123.75
$ws.Range("M4").Value = -11.75

$ws.Range("H9").Value = 1900
$ws.Range("J9").Value = 1900
$ws.Range("L9").Value = 5700
$ws.Range("N9").Value = -6148

$ws.Range("H10").Value = 78.5
$ws.Range("I10").Value = 78.5
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 235.5
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -96.5
$ws.Range("N10").ClearContents()

$ws.Range("H15").Value = 583.1111
$ws.Range("I15").Value = 34
$ws.Range("J15").Value = 1022.4
$ws.Range("K15").Value = 102
$ws.Range("L15").Value = 3067.2
$ws.Range("M15").Value = 38
$ws.Range("N15").Value = -3347.2

$ws.Range("H16").Value = 1749.5
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1749.5
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 5248.5
$ws.Range("N16").Value = -5594.5
$ws.Range("M16").ClearContents()

$ws.Range("H21").Value = 850
$ws.Range("I21").Value = 800
$ws.Range("J21").Value = 900
$ws.Range("K21").Value = 2400
$ws.Range("L21").Value = 2700
$ws.Range("M21").Value = -2227
$ws.Range("N21").Value = -3046

$ws.Range("H22").Value = 1750
$ws.Range("I22").Value = 950
$ws.Range("J22").Value = 2550
$ws.Range("K22").Value = 2850
$ws.Range("L22").Value = 7650
$ws.Range("M22").Value = -2681
$ws.Range("N22").Value = -7988

$ws.Range("H27").Value = 1750
$ws.Range("I27").Value = 950
$ws.Range("J27").Value = 2550
$ws.Range("K27").Value = 2850
$ws.Range("L27").Value = 7650
$ws.Range("M27").Value = -2748
$ws.Range("N27").Value = -7854

$ws.Range("H87").Value = 4115.778
$ws.Range("I87").Value = 3505.25
$ws.Range("J87").Value = 9000
$ws.Range("K87").Value = 10515.75
$ws.Range("L87").Value = 27000
$ws.Range("M87").Value = -9267.75
$ws.Range("N87").Value = -29496

$ws.Range("H90").Value = 4115.778
$ws.Range("I90").Value = 3505.25
$ws.Range("J90").Value = 9000
$ws.Range("K90").Value = 31547.25
$ws.Range("L90").Value = 81000
$ws.Range("M90").Value = -25307.25
$ws.Range("N90").Value = -93480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1000
$ws.Range("I97").Value = 1000
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1000
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -504
$ws.Range("N97").ClearContents()

$ws.Range("H113").Value = 1406.4546
$ws.Range("I113").Value = 1377.1
$ws.Range("K113").Value = 1377.1
$ws.Range("M113").Value = 792.9000000000001

$ws.Range("H126").Value = 111113060
$ws.Range("I126").Value = 333333340
$ws.Range("J126").Value = 2907
$ws.Range("K126").Value = 1000000020
$ws.Range("L126").Value = 8721
$ws.Range("M126").Value = -999997550
$ws.Range("N126").Value = -13661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 15000
$ws.Range("J18").Value = 15000
$ws.Range("L18").Value = 15000
$ws.Range("N18").Value = -15344

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5018.6943
$ws.Range("I136").Value = 2313.3794
$ws.Range("J136").Value = 16226.429
$ws.Range("K136").Value = 6940.138199999999
$ws.Range("L136").Value = 48679.287
$ws.Range("M136").Value = -4390.138199999999
$ws.Range("N136").Value = -53779.287
